# Updating the HED "bad defs" test sheet: add a new data row (row 4) that
# exercises an invalid HED "Def/" tag, so validators have something bad to
# flag when they're run as a chained series of validators.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: same event/description as row 3 ("PerturbRight"), but with a
# new event code and a HED tag that uses an invalid Def/ construct.
$ws.Cells.Item(4, 1).Value = 254
$ws.Cells.Item(4, 2).Value = "PerturbRight"
$ws.Cells.Item(4, 3).Value = "Vehicle undergoes a perturbation to right."
$ws.Cells.Item(4, 4).Value = "Def/DefInvalid"
$ws.Rows.Item(4).RowHeight = 14.9

# The old row 4 was just a blank, tall spacer row and row 5 was a blank
# trailing row - remove the old trailing blank row now that row 4 holds data.
$ws.Rows.Item(5).Delete()

# Touch the very last row of the sheet (in column D, within the sheet's
# existing data columns) so the worksheet's extent/dimension stretches all
# the way down to row 1048576, then clear the value back out so the row
# stays empty while keeping that extended extent, and give it its own
# (short) row height.
$lastCell = $ws.Cells.Item(1048576, 4)
$lastCell.Value = "x"
$lastCell.Value = $null
$ws.Rows.Item(1048576).RowHeight = 12.8

# Leave the selection on the newly added cell, scrolled back to the top of
# the sheet.
$ws.Range("A1").Select()
$ws.Range("D4").Select()
